$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44369
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 8000
$ws.Range("P2").Value = 800

# Row 3
$ws.Range("D3").Value = 44473
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 11000
$ws.Range("M3").Value = 11000
$ws.Range("P3").Value = 1100

# Row 4
$ws.Range("D4").Value = 44469
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("P4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44463
$ws.Range("J5").Value = 25

# Row 7
$ws.Range("D7").Value = 44348
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 10000
$ws.Range("P7").Value = 1000
